$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 261.85715
$ws.Range("I11").Value = 261.85715
$ws.Range("K11").Value = 261.85715
$ws.Range("M11").Value = -121.85715
$ws.Range("H13").Value = 499.53845
$ws.Range("I13").Value = 495
$ws.Range("J13").Value = 514.6667
$ws.Range("K13").Value = 495
$ws.Range("L13").Value = 514.6667
$ws.Range("M13").Value = -326
$ws.Range("N13").Value = -852.6667
$ws.Range("H34").Value = 915
$ws.Range("I34").Value = 915
$ws.Range("K34").Value = 915
$ws.Range("M34").Value = -712
$ws.Range("H36").Value = 915
$ws.Range("I36").Value = 915
$ws.Range("K36").Value = 915
$ws.Range("M36").Value = -200
$ws.Range("H76").Value = 9921.764999999999
$ws.Range("I76").Value = 13368
$ws.Range("K76").Value = 13368
$ws.Range("M76").Value = -13053
$ws.Range("H79").Value = 9921.764999999999
$ws.Range("I79").Value = 13368
$ws.Range("K79").Value = 13368
$ws.Range("M79").Value = -12276
$ws.Range("H129").Value = 4600
$ws.Range("J129").Value = 8000
$ws.Range("L129").Value = 24000
$ws.Range("N129").Value = -34000
$ws.Range("H132").Value = 2225.05
$ws.Range("I132").Value = 1747.9143
$ws.Range("K132").Value = 5243.742899999999
$ws.Range("M132").Value = -2713.742899999999
$ws.Range("H133").Value = 119999.5
$ws.Range("J133").Value = 119999.5
$ws.Range("L133").Value = 119999.5
$ws.Range("N133").Value = -130119.5
$ws.Range("H137").Value = 2685.2917
$ws.Range("I137").Value = 3586.0908
$ws.Range("K137").Value = 10758.2724
$ws.Range("M137").Value = -8208.2724
$ws.Range("H138").Value = 9808662
$ws.Range("I138").Value = 1343.2307
$ws.Range("J138").Value = 15879859
$ws.Range("K138").Value = 4029.6921
$ws.Range("L138").Value = 47639577
$ws.Range("M138").Value = 1110.3079
$ws.Range("N138").Value = -47649857
$ws.Range("H141").Value = 2409
$ws.Range("I141").Value = 2409
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 7227
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -2047
$ws.Range("N141").ClearContents()

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2014.2
$ws.Range("I45").Value = 1844.5834
$ws.Range("K45").Value = 1844.5834
$ws.Range("M45").Value = -1467.5834
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H61").Value = 35723176
$ws.Range("I61").Value = 50008490
$ws.Range("K61").Value = 50008490
$ws.Range("M61").Value = -50008278
$ws.Range("H74").Value = 24418352
$ws.Range("I74").Value = 24418352
$ws.Range("K74").Value = 24418352
$ws.Range("M74").Value = -24417478
$ws.Range("H77").Value = 24418352
$ws.Range("I77").Value = 24418352
$ws.Range("K77").Value = 122091760
$ws.Range("M77").Value = -122087392
$ws.Range("H88").Value = 9948.583000000001
$ws.Range("I88").Value = 21100.8
$ws.Range("J88").Value = 1982.7142
$ws.Range("K88").Value = 21100.8
$ws.Range("L88").Value = 1982.7142
$ws.Range("M88").Value = -20694.8
$ws.Range("N88").Value = -2794.7142
$ws.Range("H91").Value = 9948.583000000001
$ws.Range("I91").Value = 21100.8
$ws.Range("J91").Value = 1982.7142
$ws.Range("K91").Value = 21100.8
$ws.Range("L91").Value = 1982.7142
$ws.Range("M91").Value = -19696.8
$ws.Range("N91").Value = -4790.7142
$ws.Range("H132").Value = 111270550
$ws.Range("I132").Value = 4990.8
$ws.Range("K132").Value = 14972.4
$ws.Range("M132").Value = -12442.4
$ws.Range("H136").Value = 35723176
$ws.Range("I136").Value = 50008490
$ws.Range("K136").Value = 150025470
$ws.Range("M136").Value = -150022920

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 24950
$ws.Range("I44").Value = 24950
$ws.Range("K44").Value = 24950
$ws.Range("M44").Value = -24453
$ws.Range("H86").Value = 16045.35
$ws.Range("I86").Value = 7332.6875
$ws.Range("J86").Value = 50896
$ws.Range("K86").Value = 7332.6875
$ws.Range("L86").Value = 50896
$ws.Range("M86").Value = -6209.6875
$ws.Range("N86").Value = -53142
$ws.Range("H89").Value = 16045.35
$ws.Range("I89").Value = 7332.6875
$ws.Range("J89").Value = 50896
$ws.Range("K89").Value = 36663.4375
$ws.Range("L89").Value = 254480
$ws.Range("M89").Value = -31047.4375
$ws.Range("N89").Value = -265712
$ws.Range("H134").Value = 4115.615
$ws.Range("I134").Value = 3900.4
$ws.Range("K134").Value = 11701.2
$ws.Range("M134").Value = -9166.200000000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2144.7144
$ws.Range("I58").Value = 1603.7273
$ws.Range("K58").Value = 1603.7273
$ws.Range("M58").Value = -1400.7273
$ws.Range("H103").Value = 60000
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("H136").Value = 2144.7144
$ws.Range("I136").Value = 1603.7273
$ws.Range("K136").Value = 4811.1819
$ws.Range("M136").Value = -2261.1819

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 360
$ws.Range("I35").Value = 360
$ws.Range("K35").Value = 1080
$ws.Range("M35").Value = -792
$ws.Range("H97").Value = 259.5
$ws.Range("J97").Value = 248
$ws.Range("L97").Value = 744
$ws.Range("N97").Value = -1736
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H129").Value = 4043.8928
$ws.Range("I129").Value = 4437
$ws.Range("J129").Value = 3886.65
$ws.Range("K129").Value = 13311
$ws.Range("L129").Value = 11659.95
$ws.Range("M129").Value = -8311
$ws.Range("N129").Value = -21659.95
$ws.Range("H131").Value = 24001.52
$ws.Range("I131").Value = 61983.176
$ws.Range("J131").Value = 4435.212
$ws.Range("K131").Value = 185949.528
$ws.Range("L131").Value = 13305.636
$ws.Range("M131").Value = -180909.528
$ws.Range("N131").Value = -23385.636
$ws.Range("H132").Value = 1697.6666
$ws.Range("I132").Value = 1697.6666
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15278.9994
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12748.9994
$ws.Range("N132").ClearContents()
$ws.Range("H139").Value = 3969.8
$ws.Range("J139").Value = 4000
$ws.Range("L139").Value = 12000
$ws.Range("N139").Value = -22280

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 3125319.5
$ws.Range("J2").Value = 499.25
$ws.Range("L2").Value = 499.25
$ws.Range("N2").Value = -725.25
$ws.Range("H80").Value = 3706.6956
$ws.Range("I80").Value = 3663.375
$ws.Range("K80").Value = 3663.375
$ws.Range("M80").Value = -2665.375
$ws.Range("H83").Value = 3706.6956
$ws.Range("I83").Value = 3663.375
$ws.Range("K83").Value = 18316.875
$ws.Range("M83").Value = -13324.875
$ws.Range("H113").Value = 3005.8276
$ws.Range("I113").Value = 1579
$ws.Range("K113").Value = 1579
$ws.Range("M113").Value = 591
$ws.Range("H126").Value = 8584.727999999999
$ws.Range("I126").Value = 19722
$ws.Range("K126").Value = 59166
$ws.Range("M126").Value = -56696
$ws.Range("H132").Value = 10882.4
$ws.Range("I132").Value = 15970.667
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 47912.001
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = -45382.001
$ws.Range("N132").Value = -14810

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 49236.5
$ws.Range("J6").Value = 49236.5
$ws.Range("L6").Value = 49236.5
$ws.Range("N6").Value = -49460.5
$ws.Range("H22").Value = 3362.8333
$ws.Range("J22").Value = 3540.4
$ws.Range("L22").Value = 3540.4
$ws.Range("N22").Value = -4130.4
$ws.Range("H27").Value = 3362.8333
$ws.Range("J27").Value = 3540.4
$ws.Range("L27").Value = 3540.4
$ws.Range("N27").Value = -3754.4
$ws.Range("H55").Value = 715.7368
$ws.Range("I55").Value = 433.55554
$ws.Range("K55").Value = 433.55554
$ws.Range("M55").Value = -260.55554
$ws.Range("H61").Value = 3085.2188
$ws.Range("I61").Value = 2459.5264
$ws.Range("K61").Value = 2459.5264
$ws.Range("M61").Value = -2257.5264
$ws.Range("H100").Value = 4078.2
$ws.Range("I100").Value = 3397.6667
$ws.Range("K100").Value = 3397.6667
$ws.Range("M100").Value = -2856.6667
$ws.Range("H113").Value = 3085.2188
$ws.Range("I113").Value = 2459.5264
$ws.Range("K113").Value = 2459.5264
$ws.Range("M113").Value = -289.5264000000002
$ws.Range("H122").Value = 4213.095
$ws.Range("I122").Value = 3350.3076
$ws.Range("J122").Value = 5615.125
$ws.Range("K122").Value = 10050.9228
$ws.Range("L122").Value = 16845.375
$ws.Range("M122").Value = -7600.9228
$ws.Range("N122").Value = -21745.375
$ws.Range("H136").Value = 2387.0417
$ws.Range("I136").Value = 1758.641
$ws.Range("K136").Value = 5275.923000000001
$ws.Range("M136").Value = -2725.923000000001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1153.2727
$ws.Range("J113").Value = 1355.1428
$ws.Range("L113").Value = 4065.4284
$ws.Range("N113").Value = -8405.428400000001
$ws.Range("H126").Value = 7691.8
$ws.Range("I126").Value = 8364.75
$ws.Range("K126").Value = 25094.25
$ws.Range("M126").Value = -22624.25
$ws.Range("H129").Value = 68499
$ws.Range("J129").Value = 68499
$ws.Range("L129").Value = 68499
$ws.Range("N129").Value = -78499
$ws.Range("H132").Value = 3467.2083
$ws.Range("I132").Value = 3268.075
$ws.Range("K132").Value = 9804.224999999999
$ws.Range("M132").Value = -7274.224999999999
$ws.Range("H136").Value = 1029.279
$ws.Range("I136").Value = 844.4286
$ws.Range("K136").Value = 2533.2858
$ws.Range("M136").Value = 16.71420000000035
